# Auto update: 2025-12-05 17:31:00
# Refresh the daily 방산(defense) stock analysis data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: HANWHA AEROSPACE / 012450.KS
$ws.Range("B2").Value = "HANWHA AEROSPACE"
$ws.Range("C2").Value = "012450.KS"
$ws.Range("D2").Value = 895000
$ws.Range("E2").Value = 37
$ws.Range("F2").Value = 5.05
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 56
$ws.Range("I2").Value = 76
$ws.Range("J2").Value = 66
$ws.Range("K2").Value = 55.8
$ws.Range("N2").Value = 54.77309453746771

# Row 3: HYUNDAI ROTEM / 064350.KS
$ws.Range("B3").Value = "HYUNDAI ROTEM"
$ws.Range("C3").Value = "064350.KS"
$ws.Range("D3").Value = 182000
$ws.Range("E3").Value = 36.5
$ws.Range("F3").Value = 3.59
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 60
$ws.Range("I3").Value = 66
$ws.Range("J3").Value = 63
$ws.Range("K3").Value = 51.8
$ws.Range("N3").Value = 54.77309453746771

# Row 4: KOREA AEROSPACE / 047810.KS
$ws.Range("B4").Value = "KOREA AEROSPACE"
$ws.Range("C4").Value = "047810.KS"
$ws.Range("D4").Value = 106200
$ws.Range("E4").Value = 40.6
$ws.Range("F4").Value = -2.48
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 43
$ws.Range("I4").Value = 63
$ws.Range("J4").Value = 66
$ws.Range("K4").Value = 47.6
$ws.Range("N4").Value = 54.77309453746771

# Row 5: HANWHA SYSTEMS / 272210.KS
$ws.Range("B5").Value = "HANWHA SYSTEMS"
$ws.Range("C5").Value = "272210.KS"
$ws.Range("D5").Value = 47200
$ws.Range("E5").Value = 25.6
$ws.Range("F5").Value = 2.16
$ws.Range("G5").Value = 20
$ws.Range("H5").Value = 53
$ws.Range("I5").Value = 60
$ws.Range("J5").Value = 43
$ws.Range("K5").Value = 46.4
$ws.Range("N5").Value = 54.77309453746771

# Row 6: LIG Nex1 / 079550.KS
$ws.Range("B6").Value = "LIG Nex1"
$ws.Range("C6").Value = "079550.KS"
$ws.Range("D6").Value = 372500
$ws.Range("E6").Value = 29.6
$ws.Range("F6").Value = -2.74
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 56
$ws.Range("I6").Value = 46
$ws.Range("J6").Value = 46
$ws.Range("K6").Value = 37.8
$ws.Range("N6").Value = 54.77309453746771

$wb.Save()
